$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5941326565173313
$ws.Range("C2").Value = 0.2434560285435161
$ws.Range("D2").Value = 0.04876895091653211
$ws.Range("F2").Value = 1.003736934346527
$ws.Range("G2").Value = 0.002467271801632549
$ws.Range("K2").Value = 0.2772281472357463
$ws.Range("L2").Value = 0.2860582655040389
$ws.Range("M2").Value = 0.1971963222251176
$ws.Range("N2").Value = 2.021994375149569
$ws.Range("O2").Value = 3.597219451174766
$ws.Range("B3").Value = 0.5580601044114815
$ws.Range("C3").Value = 0.2437574110681808
$ws.Range("D3").Value = 0.04677760021012034
$ws.Range("F3").Value = 1.00340876471985
$ws.Range("G3").Value = 0.002469647624118693
$ws.Range("K3").Value = 0.2453860685970994
$ws.Range("L3").Value = 0.282656447314217
$ws.Range("M3").Value = 0.1899200115138839
$ws.Range("N3").Value = 2.041163419150831
$ws.Range("O3").Value = 3.609649069998426
$ws.Range("B4").Value = 0.5361364545110803
$ws.Range("C4").Value = 0.2439638758301292
$ws.Range("D4").Value = 0.04553996981641717
$ws.Range("F4").Value = 1.003726203126639
$ws.Range("G4").Value = 0.002471185562136189
$ws.Range("K4").Value = 0.2258452604874179
$ws.Range("L4").Value = 0.2807098957469378
$ws.Range("M4").Value = 0.1855438203428221
$ws.Range("N4").Value = 2.053538382658555
$ws.Range("O4").Value = 3.618989955416666
$ws.Range("B5").Value = 0.5272595504620483
$ws.Range("C5").Value = 0.2440534203269245
$ws.Range("D5").Value = 0.04503189011922615
$ws.Range("F5").Value = 1.003986160025725
$ws.Range("G5").Value = 0.002471832251263228
$ws.Range("K5").Value = 0.2178852355198302
$ws.Range("L5").Value = 0.2799524877371482
$ws.Range("M5").Value = 0.1837836144923664
$ws.Range("N5").Value = 2.058733461910812
$ws.Range("O5").Value = 3.62322646236349
$ws.Range("B6").Value = 0.5257890165565016
$ws.Range("C6").Value = 0.2440686164741557
$ws.Range("D6").Value = 0.04494729885248461
$ws.Range("F6").Value = 1.004037216646545
$ws.Range("G6").Value = 0.002471940841290463
$ws.Range("K6").Value = 0.2165636734468563
$ws.Range("L6").Value = 0.2798288868355172
$ws.Range("M6").Value = 0.1834927340593993
$ws.Range("N6").Value = 2.059605293018917
$ws.Range("O6").Value = 3.623955911137941
$ws.Range("B7").Value = 0.5360165052754269
$ws.Range("C7").Value = 0.2439650615278168
$ws.Range("D7").Value = 0.0455331327736701
$ws.Range("F7").Value = 1.003729180067168
$ws.Range("G7").Value = 0.002471194202657997
$ws.Range("K7").Value = 0.2257378960655387
$ws.Range("L7").Value = 0.2806995359128308
$ws.Range("M7").Value = 0.185519987787945
$ws.Range("N7").Value = 2.053607829104944
$ws.Range("O7").Value = 3.619045348944326
$ws.Range("B8").Value = 0.5816484642720923
$ws.Range("C8").Value = 0.2435555146725505
$ws.Range("D8").Value = 0.04808544388711766
$ws.Range("F8").Value = 1.003516138108779
$ws.Range("G8").Value = 0.002468074587996768
$ws.Range("K8").Value = 0.2662470844843767
$ws.Range("L8").Value = 0.2848558530717966
$ws.Range("M8").Value = 0.1946685308208949
$ws.Range("N8").Value = 2.028478291667076
$ws.Range("O8").Value = 3.601150624780047
$ws.Range("B9").Value = 0.6728979172624179
$ws.Range("C9").Value = 0.242921291756506
$ws.Range("D9").Value = 0.05297136549076953
$ws.Range("F9").Value = 1.007212432688014
$ws.Range("G9").Value = 0.002462582579860317
$ws.Range("K9").Value = 0.3457535561194902
$ws.Range("L9").Value = 0.2941319209171098
$ws.Range("M9").Value = 0.2133305255093276
$ws.Range("N9").Value = 1.983998589243514
$ws.Range("O9").Value = 3.579610591232409
$ws.Range("B10").Value = 0.7409949578092494
$ws.Range("C10").Value = 0.2425569454950214
$ws.Range("D10").Value = 0.05648781641804135
$ws.Range("F10").Value = 1.012433383141229
$ws.Range("G10").Value = 0.002458925229948214
$ws.Range("K10").Value = 0.4041951918953259
$ws.Range("L10").Value = 0.3016308351382833
$ws.Range("M10").Value = 0.2274774003478441
$ws.Range("N10").Value = 1.954241408580979
$ws.Range("O10").Value = 3.572038141485422
$ws.Range("B11").Value = 0.772199530679103
$ws.Range("C11").Value = 0.242412978434956
$ws.Range("D11").Value = 0.05807152786079683
$ws.Range("F11").Value = 1.015352064051839
$ws.Range("G11").Value = 0.002457342614072717
$ws.Range("K11").Value = 0.4307854483213589
$ws.Range("L11").Value = 0.3051903256867092
$ws.Range("M11").Value = 0.234007061882842
$ws.Range("N11").Value = 1.941337805208827
$ws.Range("O11").Value = 3.570383833278669
$ws.Range("B12").Value = 0.7840480284219211
$ws.Range("C12").Value = 0.2423615700673878
$ws.Range("D12").Value = 0.05866892928939649
$ws.Range("F12").Value = 1.016535392748892
$ws.Range("G12").Value = 0.002456754925108576
$ws.Range("K12").Value = 0.440854842716476
$ws.Range("L12").Value = 0.3065594620633334
$ws.Range("M12").Value = 0.2364931160534098
$ws.Range("N12").Value = 1.9365425638
$ws.Range("O12").Value = 3.570014670126966
$ws.Range("B13").Value = 0.7814948290127006
$ws.Range("C13").Value = 0.2423725038691984
$ws.Range("D13").Value = 0.05854037157197922
$ws.Range("F13").Value = 1.016277070181388
$ws.Range("G13").Value = 0.002456880978778179
$ws.Range("K13").Value = 0.4386862150332433
$ws.Range("L13").Value = 0.3062636508481518
$ws.Range("M13").Value = 0.235957105322889
$ws.Range("N13").Value = 1.937571253732655
$ws.Range("O13").Value = 3.570082735880845
$ws.Range("B14").Value = 0.7731736758989598
$ws.Range("C14").Value = 0.2424086868600952
$ws.Range("D14").Value = 0.05812072298196824
$ws.Range("F14").Value = 1.015447852730048
$ws.Range("G14").Value = 0.002457294032050283
$ws.Range("K14").Value = 0.4316138612467455
$ws.Range("L14").Value = 0.3053025402242753
$ws.Range("M14").Value = 0.2342113229810181
$ws.Range("N14").Value = 1.940941472355524
$ws.Range("O14").Value = 3.570348306584691
$ws.Range("B15").Value = 0.7680808786072078
$ws.Range("C15").Value = 0.2424312541835292
$ws.Range("D15").Value = 0.05786337382570395
$ws.Range("F15").Value = 1.014950099787427
$ws.Range("G15").Value = 0.00245754854993677
$ws.Range("K15").Value = 0.4272818581715114
$ws.Range("L15").Value = 0.3047165953212527
$ws.Range("M15").Value = 0.2331437238151111
$ws.Range("N15").Value = 1.94301769009744
$ws.Range("O15").Value = 3.570544477525601
$ws.Range("B16").Value = 0.7389601779157999
$ws.Range("C16").Value = 0.2425667903394597
$ws.Range("D16").Value = 0.05638399438708319
$ws.Range("F16").Value = 1.012253572503383
$ws.Range("G16").Value = 0.002459030285637588
$ws.Range("K16").Value = 0.4024575148335714
$ws.Range("L16").Value = 0.301401189673399
$ws.Range("M16").Value = 0.2270525556331791
$ws.Range("N16").Value = 1.955097431759008
$ws.Range("O16").Value = 3.572182267618757
$ws.Range("B17").Value = 0.7211532396890732
$ws.Range("C17").Value = 0.2426554986312475
$ws.Range("D17").Value = 0.05547234370015275
$ws.Range("F17").Value = 1.010738518458922
$ws.Range("G17").Value = 0.002459960024239871
$ws.Range("K17").Value = 0.3872295042334599
$ws.Range("L17").Value = 0.2994051985134547
$ws.Range("M17").Value = 0.2233398502459423
$ws.Range("N17").Value = 1.962670135631837
$ws.Range("O17").Value = 3.573645448120629
$ws.Range("B18").Value = 0.7109325528277282
$ws.Range("C18").Value = 0.2427085722808755
$ws.Range("D18").Value = 0.05494648743151487
$ws.Range("F18").Value = 1.009918275116483
$ws.Range("G18").Value = 0.002460502424964872
$ws.Range("K18").Value = 0.3784712404702475
$ws.Range("L18").Value = 0.2982711107436131
$ws.Range("M18").Value = 0.2212132724728733
$ws.Range("N18").Value = 1.96708534146779
$ws.Range("O18").Value = 3.57465556574428
$ws.Range("B19").Value = 0.7074756971111071
$ws.Range("C19").Value = 0.2427268950530745
$ws.Range("D19").Value = 0.05476818481406553
$ws.Range("F19").Value = 1.009649346951157
$ws.Range("G19").Value = 0.00246068738605604
$ws.Range("K19").Value = 0.3755059401174492
$ws.Range("L19").Value = 0.2978895269564532
$ws.Range("M19").Value = 0.2204947775227808
$ws.Range("N19").Value = 1.968590487350214
$ws.Range("O19").Value = 3.575026525587361
$ws.Range("B20").Value = 0.723046609074288
$ws.Range("C20").Value = 0.2426458433768772
$ws.Range("D20").Value = 0.05556954577389206
$ws.Range("F20").Value = 1.010894502691677
$ws.Range("G20").Value = 0.002459860261858275
$ws.Range("K20").Value = 0.3888505050599349
$ws.Range("L20").Value = 0.29961623158907
$ws.Range("M20").Value = 0.2237341567737161
$ws.Range("N20").Value = 1.961857841106308
$ws.Range("O20").Value = 3.573472248593617
$ws.Range("B21").Value = 0.7756169366874985
$ws.Range("C21").Value = 0.2423979748445078
$ws.Range("D21").Value = 0.05824404698216057
$ws.Range("F21").Value = 1.015689295424963
$ws.Range("G21").Value = 0.002457172393439739
$ws.Range("K21").Value = 0.4336911800928362
$ws.Range("L21").Value = 0.3055842660196504
$ws.Range("M21").Value = 0.2347237385804135
$ws.Range("N21").Value = 1.939949085452444
$ws.Range("O21").Value = 3.570263320721381
$ws.Range("B22").Value = 0.8101609342707548
$ws.Range("C22").Value = 0.2422540887529649
$ws.Range("D22").Value = 0.05997847929962319
$ws.Range("F22").Value = 1.019278082568846
$ws.Range("G22").Value = 0.002455483385237434
$ws.Range("K22").Value = 0.4629983472077015
$ws.Range("L22").Value = 0.3096084572796798
$ws.Range("M22").Value = 0.2419841868154862
$ws.Range("N22").Value = 1.926161309605433
$ws.Range("O22").Value = 3.569665687563543
$ws.Range("B23").Value = 0.7917073168787567
$ws.Range("C23").Value = 0.2423292337669807
$ws.Range("D23").Value = 0.0590540240121058
$ws.Range("F23").Value = 1.01732106566844
$ws.Range("G23").Value = 0.00245637866628523
$ws.Range("K23").Value = 0.4473566028313485
$ws.Range("L23").Value = 0.3074493740249551
$ws.Range("M23").Value = 0.2381020437290857
$ws.Range("N23").Value = 1.933471518569725
$ws.Range("O23").Value = 3.569847501319515
$ws.Range("B24").Value = 0.7221905645259881
$ws.Range("C24").Value = 0.2426502020557244
$ws.Range("D24").Value = 0.05552560612206037
$ws.Range("F24").Value = 1.010823824032805
$ws.Range("G24").Value = 0.002459905339875358
$ws.Range("K24").Value = 0.3881176613849675
$ws.Range("L24").Value = 0.2995207816849899
$ws.Range("M24").Value = 0.2235558661567865
$ws.Range("N24").Value = 1.96222488771325
$ws.Range("O24").Value = 3.573550025922543
$ws.Range("B25").Value = 0.6480256098870143
$ws.Range("C25").Value = 0.2430749241266561
$ws.Range("D25").Value = 0.05166241177751374
$ws.Range("F25").Value = 1.005772331589291
$ws.Range("G25").Value = 0.002464001730925614
$ws.Range("K25").Value = 0.3242390304662308
$ws.Range("L25").Value = 0.2915022461205581
$ws.Range("M25").Value = 0.2082051159760994
$ws.Range("N25").Value = 1.99551812712056
$ws.Range("O25").Value = 3.583987825270697
